# TS Pada Padam 1.1 to 1.8 final
# Applies two formatting/text fixes to the document:
#  1. The "(ignore those which are already incorporated ...)" paragraph
#     (2nd paragraph in the document) switches its run/paragraph-mark
#     formatting from bCs/sz32/szCs32 to b/szCs24.
#  2. The "TS Pada Paatam - TS 1.5 ... 31st March 2020." heading merges the
#     "- TS " and "1" runs into a single run "- TS 1" and drops the
#     gramStart/gramEnd proofing-error markers that bracketed "1 ... Tamil".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "(ignore those which are already incorporated ...)" paragraph
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(2)
$r1 = $p1.Range
if ($r1.Text -notlike "*already incorporated*") {
    throw "Change 1: paragraph 2 does not contain the expected text; aborting."
}

$change1Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6892D140" w14:textId="77777777" w:rsidR="00213E19" w:rsidRPr="00113311" w:rsidRDefault="00213E19" w:rsidP="00213E19"><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00113311"><w:rPr><w:b/><w:szCs w:val="24"/></w:rPr><w:t>(ignore those which are already incorporated in your book’s version and date)</w:t></w:r></w:p>'
$res1 = $r1.InsertXML($change1Xml)
Write-Output ("Change 1 applied (InsertXML result: " + $res1 + ")")

# ---------------------------------------------------------------------
# Change 2: "TS Pada Paatam - TS 1.5 ... Observed till 31st March 2020."
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(221)
$r2 = $p2.Range
if ($r2.Text -notlike "*31st*March*2020*") {
    throw "Change 2: paragraph 221 does not contain the expected text; aborting."
}

$change2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="38809A0D" w14:textId="33475AFF" w:rsidR="0000227C" w:rsidRPr="00D40DD6" w:rsidRDefault="00DD72F2" w:rsidP="0000227C"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>=============</w:t></w:r><w:r w:rsidR="004E083D"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:br w:type="page"/></w:r><w:r w:rsidR="0000227C"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">TS Pada Paatam </w:t></w:r><w:r w:rsidR="0000227C" w:rsidRPr="00D40DD6"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>– TS 1</w:t></w:r><w:r w:rsidR="0000227C"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">.5 </w:t></w:r><w:r w:rsidR="0000227C" w:rsidRPr="00D40DD6"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0000227C"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>Tamil</w:t></w:r><w:r w:rsidR="0000227C" w:rsidRPr="00D40DD6"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> Corrections –</w:t></w:r><w:r w:rsidR="0000227C"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0000227C" w:rsidRPr="00D40DD6"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Observed till </w:t></w:r><w:r w:rsidR="0000227C"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>31st</w:t></w:r><w:r w:rsidR="0000227C" w:rsidRPr="00D40DD6"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0000227C"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">March </w:t></w:r><w:r w:rsidR="0000227C" w:rsidRPr="00D40DD6"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>20</w:t></w:r><w:r w:rsidR="0000227C"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>20.</w:t></w:r></w:p>'
$res2 = $r2.InsertXML($change2Xml)
Write-Output ("Change 2 applied (InsertXML result: " + $res2 + ")")
